$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 1.524029
$ws.Range("H2").Value = 3.048058
$ws.Range("I2").Value = 0.09030204154573296
$ws.Range("J2").Value = 0.06866669168778029
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.704614
$ws.Range("N2").Value = 5.113842
$ws.Range("O2").Value = 0.1078894737559977
$ws.Range("P2").Value = 0.1126857637256889
$ws.Range("Q2").Value = 2.597881169806
$ws.Range("R2").Value = 15.587287018836
$ws.Range("S2").Value = 0.009742639741461369
$ws.Range("T2").Value = 0.007737758595353935

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1.524029
$ws.Range("H3").Value = 3.048058
$ws.Range("I3").Value = 0.09030204154573296
$ws.Range("J3").Value = 0.06866669168778029
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.267396666666667
$ws.Range("N3").Value = 24.80219
$ws.Range("O3").Value = 0.5232651355079543
$ws.Range("P3").Value = 0.5465271946649201
$ws.Range("Q3").Value = 12.59975227450333
$ws.Range("R3").Value = 75.59851364702
$ws.Range("S3").Value = 0.04725191000607287
$ws.Range("T3").Value = 0.03752821437504355

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 1.524029
$ws.Range("H4").Value = 3.048058
$ws.Range("I4").Value = 0.09030204154573296
$ws.Range("J4").Value = 0.06866669168778029
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.810163
$ws.Range("N4").Value = 11.430489
$ws.Range("O4").Value = 0.2411551711968653
$ws.Range("P4").Value = 0.2518758660754646
$ws.Range("Q4").Value = 5.806798906727
$ws.Range("R4").Value = 34.840793440362
$ws.Range("S4").Value = 0.02177680428838768
$ws.Range("T4").Value = 0.01729548243939657

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 1.524029
$ws.Range("H5").Value = 3.048058
$ws.Range("I5").Value = 0.09030204154573296
$ws.Range("J5").Value = 0.06866669168778029
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 2.0174585
$ws.Range("N5").Value = 4.034917
$ws.Range("O5").Value = 0.1276902195391827
$ws.Range("P5").Value = 0.08891117553392644
$ws.Range("Q5").Value = 3.0746652602965
$ws.Range("R5").Value = 12.298661041186
$ws.Range("S5").Value = 0.01153068750981104
$ws.Range("T5").Value = 0.006105236277986241

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.502875000000001
$ws.Range("H6").Value = 25.508625
$ws.Range("I6").Value = 0.5038138851085998
$ws.Range("J6").Value = 0.5746586476550659
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.704614
$ws.Range("N6").Value = 5.113842
$ws.Range("O6").Value = 0.1078894737559977
$ws.Range("P6").Value = 0.1126857637256889
$ws.Range("Q6").Value = 14.49411976525
$ws.Range("R6").Value = 130.44707788725
$ws.Range("S6").Value = 0.05435621493533152
$ws.Range("T6").Value = 0.06475584859258265

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.502875000000001
$ws.Range("H7").Value = 25.508625
$ws.Range("I7").Value = 0.5038138851085998
$ws.Range("J7").Value = 0.5746586476550659
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.267396666666667
$ws.Range("N7").Value = 24.80219
$ws.Range("O7").Value = 0.5232651355079543
$ws.Range("P7").Value = 0.5465271946649201
$ws.Range("Q7").Value = 70.29664043208334
$ws.Range("R7").Value = 632.66976388875
$ws.Range("S7").Value = 0.2636282408621404
$ws.Range("T7").Value = 0.3140665785928599

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.502875000000001
$ws.Range("H8").Value = 25.508625
$ws.Range("I8").Value = 0.5038138851085998
$ws.Range("J8").Value = 0.5746586476550659
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.810163
$ws.Range("N8").Value = 11.430489
$ws.Range("O8").Value = 0.2411551711968653
$ws.Range("P8").Value = 0.2518758660754646
$ws.Range("Q8").Value = 32.39733971862501
$ws.Range("R8").Value = 291.576057467625
$ws.Range("S8").Value = 0.1214973237147222
$ws.Range("T8").Value = 0.144742644575875

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.502875000000001
$ws.Range("H9").Value = 25.508625
$ws.Range("I9").Value = 0.5038138851085998
$ws.Range("J9").Value = 0.5746586476550659
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 2.0174585
$ws.Range("N9").Value = 4.034917
$ws.Range("O9").Value = 0.1276902195391827
$ws.Range("P9").Value = 0.08891117553392644
$ws.Range("Q9").Value = 17.1541974431875
$ws.Range("R9").Value = 102.925184659125
$ws.Range("S9").Value = 0.0643321055964057
$ws.Range("T9").Value = 0.05109357589374835

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3890603333333333
$ws.Range("H10").Value = 1.167181
$ws.Range("I10").Value = 0.02305267313447669
$ws.Range("J10").Value = 0.02629426929239375
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.704614
$ws.Range("N10").Value = 5.113842
$ws.Range("O10").Value = 0.1078894737559977
$ws.Range("P10").Value = 0.1126857637256889
$ws.Range("Q10").Value = 0.6631976910446667
$ws.Range("R10").Value = 5.968779219402
$ws.Range("S10").Value = 0.002487140773147716
$ws.Range("T10").Value = 0.002962989816822319

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3890603333333333
$ws.Range("H11").Value = 1.167181
$ws.Range("I11").Value = 0.02305267313447669
$ws.Range("J11").Value = 0.02629426929239375
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.267396666666667
$ws.Range("N11").Value = 24.80219
$ws.Range("O11").Value = 0.5232651355079543
$ws.Range("P11").Value = 0.5465271946649201
$ws.Range("Q11").Value = 3.216516102932222
$ws.Range("R11").Value = 28.94864492639
$ws.Range("S11").Value = 0.01206266013153252
$ws.Range("T11").Value = 0.01437053323213591

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3890603333333333
$ws.Range("H12").Value = 1.167181
$ws.Range("I12").Value = 0.02305267313447669
$ws.Range("J12").Value = 0.02629426929239375
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.810163
$ws.Range("N12").Value = 11.430489
$ws.Range("O12").Value = 0.2411551711968653
$ws.Range("P12").Value = 0.2518758660754646
$ws.Range("Q12").Value = 1.482383286834333
$ws.Range("R12").Value = 13.341449581509
$ws.Range("S12").Value = 0.005559271336290104
$ws.Range("T12").Value = 0.006622891850843168

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3890603333333333
$ws.Range("H13").Value = 1.167181
$ws.Range("I13").Value = 0.02305267313447669
$ws.Range("J13").Value = 0.02629426929239375
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 2.0174585
$ws.Range("N13").Value = 4.034917
$ws.Range("O13").Value = 0.1276902195391827
$ws.Range("P13").Value = 0.08891117553392644
$ws.Range("Q13").Value = 0.7849130764961667
$ws.Range("R13").Value = 4.709478458977
$ws.Range("S13").Value = 0.002943600893506349
$ws.Range("T13").Value = 0.002337854392592352

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4482056666666667
$ws.Range("H14").Value = 1.344617
$ws.Range("I14").Value = 0.02655716310671665
$ws.Range("J14").Value = 0.0302915498908315
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.704614
$ws.Range("N14").Value = 5.113842
$ws.Range("O14").Value = 0.1078894737559977
$ws.Range("P14").Value = 0.1126857637256889
$ws.Range("Q14").Value = 0.7640176542793333
$ws.Range("R14").Value = 6.876158888513999
$ws.Range("S14").Value = 0.002865238352035856
$ws.Range("T14").Value = 0.003413426433883155

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4482056666666667
$ws.Range("H15").Value = 1.344617
$ws.Range("I15").Value = 0.02655716310671665
$ws.Range("J15").Value = 0.0302915498908315
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.267396666666667
$ws.Range("N15").Value = 24.80219
$ws.Range("O15").Value = 0.5232651355079543
$ws.Range("P15").Value = 0.5465271946649201
$ws.Range("Q15").Value = 3.705494034581111
$ws.Range("R15").Value = 33.34944631123
$ws.Range("S15").Value = 0.01389643755174293
$ws.Range("T15").Value = 0.0165551557838886

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4482056666666667
$ws.Range("H16").Value = 1.344617
$ws.Range("I16").Value = 0.02655716310671665
$ws.Range("J16").Value = 0.0302915498908315
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.810163
$ws.Range("N16").Value = 11.430489
$ws.Range("O16").Value = 0.2411551711968653
$ws.Range("P16").Value = 0.2518758660754646
$ws.Range("Q16").Value = 1.707736647523667
$ws.Range("R16").Value = 15.369629827713
$ws.Range("S16").Value = 0.00640439721550333
$ws.Range("T16").Value = 0.007629710363521328

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4482056666666667
$ws.Range("H17").Value = 1.344617
$ws.Range("I17").Value = 0.02655716310671665
$ws.Range("J17").Value = 0.0302915498908315
$ws.Range("K17").Value = 2
$ws.Range("M17").Value = 2.0174585
$ws.Range("N17").Value = 4.034917
$ws.Range("O17").Value = 0.1276902195391827
$ws.Range("P17").Value = 0.08891117553392644
$ws.Range("Q17").Value = 0.9042363319648333
$ws.Range("R17").Value = 5.425417991789
$ws.Range("S17").Value = 0.003391089987434534
$ws.Range("T17").Value = 0.002693257309538409

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1.295005
$ws.Range("H18").Value = 3.885015
$ws.Range("I18").Value = 0.07673187013628475
$ws.Range("J18").Value = 0.08752167025935917
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.704614
$ws.Range("N18").Value = 5.113842
$ws.Range("O18").Value = 0.1078894737559977
$ws.Range("P18").Value = 0.1126857637256889
$ws.Range("Q18").Value = 2.20748365307
$ws.Range("R18").Value = 19.86735287763
$ws.Range("S18").Value = 0.008278561089317318
$ws.Range("T18").Value = 0.0098624462557238

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1.295005
$ws.Range("H19").Value = 3.885015
$ws.Range("I19").Value = 0.07673187013628475
$ws.Range("J19").Value = 0.08752167025935917
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 8.267396666666667
$ws.Range("N19").Value = 24.80219
$ws.Range("O19").Value = 0.5232651355079543
$ws.Range("P19").Value = 0.5465271946649201
$ws.Range("Q19").Value = 10.70632002031667
$ws.Range("R19").Value = 96.35688018285001
$ws.Range("S19").Value = 0.04015111242464179
$ws.Range("T19").Value = 0.04783297291923573

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 1.295005
$ws.Range("H20").Value = 3.885015
$ws.Range("I20").Value = 0.07673187013628475
$ws.Range("J20").Value = 0.08752167025935917
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 3.810163
$ws.Range("N20").Value = 11.430489
$ws.Range("O20").Value = 0.2411551711968653
$ws.Range("P20").Value = 0.2518758660754646
$ws.Range("Q20").Value = 4.934180135815
$ws.Range("R20").Value = 44.407621222335
$ws.Range("S20").Value = 0.01850428727897139
$ws.Range("T20").Value = 0.02204459649694732

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 1.295005
$ws.Range("H21").Value = 3.885015
$ws.Range("I21").Value = 0.07673187013628475
$ws.Range("J21").Value = 0.08752167025935917
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 2.0174585
$ws.Range("N21").Value = 4.034917
$ws.Range("O21").Value = 0.1276902195391827
$ws.Range("P21").Value = 0.08891117553392644
$ws.Range("Q21").Value = 2.6126188447925
$ws.Range("R21").Value = 15.675713068755
$ws.Range("S21").Value = 0.00979790934335426
$ws.Range("T21").Value = 0.007781654587452312

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 4.717841
$ws.Range("H22").Value = 9.435682
$ws.Range("I22").Value = 0.2795423669681891
$ws.Range("J22").Value = 0.2125671712145694
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 1.704614
$ws.Range("N22").Value = 5.113842
$ws.Range("O22").Value = 0.1078894737559977
$ws.Range("P22").Value = 0.1126857637256889
$ws.Range("Q22").Value = 8.042097818374
$ws.Range("R22").Value = 48.252586910244
$ws.Range("S22").Value = 0.03015967886470392
$ws.Range("T22").Value = 0.02395329403132303

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 4.717841
$ws.Range("H23").Value = 9.435682
$ws.Range("I23").Value = 0.2795423669681891
$ws.Range("J23").Value = 0.2125671712145694
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 8.267396666666667
$ws.Range("N23").Value = 24.80219
$ws.Range("O23").Value = 0.5232651355079543
$ws.Range("P23").Value = 0.5465271946649201
$ws.Range("Q23").Value = 39.00426295726333
$ws.Range("R23").Value = 234.02557774358
$ws.Range("S23").Value = 0.1462747745318238
$ws.Range("T23").Value = 0.1161737397617564

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 4.717841
$ws.Range("H24").Value = 9.435682
$ws.Range("I24").Value = 0.2795423669681891
$ws.Range("J24").Value = 0.2125671712145694
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 3.810163
$ws.Range("N24").Value = 11.430489
$ws.Range("O24").Value = 0.2411551711968653
$ws.Range("P24").Value = 0.2518758660754646
$ws.Range("Q24").Value = 17.975743218083
$ws.Range("R24").Value = 107.854459308498
$ws.Range("S24").Value = 0.06741308736299059
$ws.Range("T24").Value = 0.05354054034888124

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 4.717841
$ws.Range("H25").Value = 9.435682
$ws.Range("I25").Value = 0.2795423669681891
$ws.Range("J25").Value = 0.2125671712145694
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 2.0174585
$ws.Range("N25").Value = 4.034917
$ws.Range("O25").Value = 0.1276902195391827
$ws.Range("P25").Value = 0.08891117553392644
$ws.Range("Q25").Value = 9.518048427098501
$ws.Range("R25").Value = 38.072193708394
$ws.Range("S25").Value = 0.03569482620867086
$ws.Range("T25").Value = 0.01889959707260878
